$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2
$ws.Range('D2').Value = '65.179.71'
$ws.Range('E2').Value = '  -2.08%  '

# Row 3
$ws.Range('D3').Value = '3.477.79'
$ws.Range('E3').Value = '  -1.01%  '

# Row 4
$ws.Range('D4').Value = "'1.00"
$ws.Range('E4').Value = '  +0.01%  '

# Row 5
$ws.Range('D5').Value = "'588.47"
$ws.Range('E5').Value = '  -2.82%  '

# Row 6
$ws.Range('E6').Value = '  -4.28%  '

# Row 7
$ws.Range('D7').Value = '3.478.59'
$ws.Range('E7').Value = '  -0.99%  '

# Row 8
$ws.Range('E8').Value = '  +0.07%  '

# Row 9
$ws.Range('E9').Value = '  -2.92%  '

# Row 10
$ws.Range('E10').Value = '  -5.75%  '

# Row 11
$ws.Range('E11').Value = '  -7.15%  '

# Row 12
$ws.Range('D12').Value = "'0.385"
$ws.Range('E12').Value = '  -4.47%  '

# Row 13
$ws.Range('D13').Value = '4.068.37'
$ws.Range('E13').Value = '  -0.80%  '

# Row 14
$ws.Range('E14').Value = '  -6.43%  '

# Row 15
$ws.Range('D15').Value = "'26.60"
$ws.Range('E15').Value = '  -7.16%  '

# Row 16
$ws.Range('D16').Value = '3.481.40'
$ws.Range('E16').Value = '  -1.25%  '

# Row 17
$ws.Range('E17').Value = '  -1.34%  '

# Row 18
$ws.Range('D18').Value = '65.161.19'
$ws.Range('E18').Value = '  -1.86%  '

# Row 19
$ws.Range('E19').Value = '  -8.51%  '

# Row 20
$ws.Range('E20').Value = '  -5.12%  '

# Row 21
$ws.Range('E21').Value = '  -4.37%  '

# Row 22
$ws.Range('D22').Value = "'389.34"
$ws.Range('E22').Value = '  -7.33%  '

# Row 23
$ws.Range('E23').Value = '  -5.10%  '

# Row 24
$ws.Range('B24').Value = 'LEO'
$ws.Range('C24').Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range('D24').Value = "'5.78"
$ws.Range('E24').Value = '  +0.69%  '

# Row 25
$ws.Range('B25').Value = 'Litecoin'
$ws.Range('C25').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D25').Value = "'72.58"
$ws.Range('E25').Value = '  -5.47%  '

# Row 26
$ws.Range('B26').Value = 'Dai'
$ws.Range('C26').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D26').Value = "'1.00"
$ws.Range('E26').Value = '  -0.02%  '

# Row 27
$ws.Range('B27').Value = 'WrappedeETH'
$ws.Range('C27').Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Range('D27').Value = '3.619.27'
$ws.Range('E27').Value = '  -1.05%  '

# Row 28
$ws.Range('E28').Value = '  -1.85%  '

# Row 29
$ws.Range('E29').Value = '  -0.01%  '

# Row 30
$ws.Range('D30').Value = "'7.38"
$ws.Range('E30').Value = '  -4.97%  '

# Row 31
$ws.Range('E31').Value = '  -8.39%  '

# Row 32
$ws.Range('D32').Value = "'2.21"
$ws.Range('E32').Value = '  -9.53%  '

# Row 33
$ws.Range('D33').Value = '3.496.92'
$ws.Range('E33').Value = '  -0.59%  '

# Row 34
$ws.Range('E34').Value = '  -0.03%  '

# Row 35
$ws.Range('E35').Value = '  -6.81%  '

# Row 36
$ws.Range('D36').Value = "'23.06"
$ws.Range('E36').Value = '  -4.39%  '

# Row 37
$ws.Range('D37').Value = "'171.64"
$ws.Range('E37').Value = '  -0.80%  '

# Row 38
$ws.Range('E38').Value = '  -8.94%  '

# Row 39
$ws.Range('D39').Value = "'6.81"
$ws.Range('E39').Value = '  -9.11%  '

# Row 40
$ws.Range('E40').Value = '  -8.66%  '

# Row 41
$ws.Range('D41').Value = "'4.73"
$ws.Range('E41').Value = '  -8.50%  '

# Row 42
$ws.Range('D42').Value = "'0.0778"
$ws.Range('E42').Value = '  -2.93%  '

# Row 43
$ws.Range('D43').Value = "'0.811"
$ws.Range('E43').Value = '  -4.59%  '

# Row 44
$ws.Range('D44').Value = "'42.55"
$ws.Range('E44').Value = '  -6.60%  '

# Row 45
$ws.Range('E45').Value = '  +0.02%  '

# Row 46
$ws.Range('D46').Value = "'25.13"
$ws.Range('E46').Value = '  +10.38%  '

# Row 47
$ws.Range('E47').Value = '  -11.89%  '

# Row 48
$ws.Range('B48').Value = 'ONDO'
$ws.Range('C48').Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Range('D48').Value = "'1.15"
$ws.Range('E48').Value = '  +3.87%  '

# Row 49
$ws.Range('B49').Value = 'Stacks'
$ws.Range('C49').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D49').Value = "'1.62"
$ws.Range('E49').Value = '  -8.08%  '

# Row 50
$ws.Range('E50').Value = '  -4.65%  '

# Row 51
$ws.Range('B51').Value = 'Maker'
$ws.Range('C51').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D51').Value = '2.219.00'
$ws.Range('E51').Value = '  -3.52%  '
